# Auto-generated Excel COM-interop script
# Updates "currentAveragePrice*" / "LevePrice*" / "LeveProfit*" columns (H:N)
# for a batch of Leve rows across multiple job sheets, matching a scheduled
# market-data refresh run ("chore: update Sheets via scheduled runner").
#
# For each touched row, columns H-N are fully re-derived from fresh market
# data: some cells are updated in place, some newly-empty-valued cells are
# removed (ClearContents), and some previously-empty cells now get a value.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 131: Mindful Study / Grade 5 Tincture of Mind
$ws.Range("H131").Value = 21571.428
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 21571.428
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 64714.284
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -74794.284

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust / Steel Ingot
$ws.Range("H32").Value = 3256.15
$ws.Range("I32").Value = 1459.8182
$ws.Range("J32").Value = 23015.8
$ws.Range("K32").Value = 1459.8182
$ws.Range("L32").Value = 23015.8
$ws.Range("M32").Value = -1172.8182
$ws.Range("N32").Value = -23589.8

# Row 61: Dealing with the Tough Stuff / Cobalt Ingot
$ws.Range("H61").Value = 52633676
$ws.Range("I61").Value = 62501764
$ws.Range("J61").Value = 3875.6667
$ws.Range("K61").Value = 62501764
$ws.Range("L61").Value = 3875.6667
$ws.Range("M61").Value = -62501552
$ws.Range("N61").Value = -4299.6667

# Row 97: Ore for Me / High Steel Ingot
$ws.Range("H97").Value = 577.75
$ws.Range("I97").Value = 546
$ws.Range("K97").Value = 546
$ws.Range("M97").Value = -50

# Row 105: Spoony Is the Bard / Molybdenum Armguards of Scouting
$ws.Range("H105").Value = 49990
$ws.Range("J105").Value = 49990
$ws.Range("L105").Value = 49990
$ws.Range("N105").Value = -56978

# Row 120: One Foot Forward / Dwarven Mythril Shoes of Maiming
$ws.Range("H120").Value = 79999
$ws.Range("J120").Value = 79999
$ws.Range("L120").Value = 79999
$ws.Range("N120").Value = -89675

# Row 122: Haste for High Durium / High Durium Nugget
$ws.Range("H122").Value = 3305.6667
$ws.Range("I122").Value = 1646.4546
$ws.Range("K122").Value = 4939.3638
$ws.Range("M122").Value = -2489.3638

# Row 136: Metal with Mettle / Cobalt Tungsten Ingot
$ws.Range("H136").Value = 52633676
$ws.Range("I136").Value = 62501764
$ws.Range("J136").Value = 3875.6667
$ws.Range("K136").Value = 187505292
$ws.Range("L136").Value = 11627.0001
$ws.Range("M136").Value = -187502742
$ws.Range("N136").Value = -16727.0001

# Row 141: Essays on Equipment / Ra'Kaznar Greaves of Maiming
$ws.Range("H141").Value = 82000
$ws.Range("J141").Value = 82000
$ws.Range("L141").Value = 82000
$ws.Range("N141").Value = -92360

$ws = $wb.Worksheets.Item("BSM")
# Row 20: Smelt and Dealt / Iron Ingot
$ws.Range("H20").Value = 1022.4138
$ws.Range("I20").Value = 1012.4545
$ws.Range("K20").Value = 1012.4545
$ws.Range("M20").Value = -765.4545000000001

# Row 94: High Steal / High Steel Nugget
$ws.Range("H94").Value = 2179.1538
$ws.Range("I94").Value = 2148.3914
$ws.Range("J94").Value = 2415
$ws.Range("K94").Value = 2148.3914
$ws.Range("L94").Value = 2415
$ws.Range("M94").Value = -1697.3914
$ws.Range("N94").Value = -3317

# Row 95: Crisscrossing / High Steel Kris
$ws.Range("H95").Value = 7605.4287
$ws.Range("J95").Value = 7605.4287
$ws.Range("L95").Value = 7605.4287
$ws.Range("N95").Value = -13097.4287

$ws = $wb.Worksheets.Item("CRP")
# Row 122: Timber of Tenkonto / Horse Chestnut Lumber
$ws.Range("H122").Value = 2667.7083
$ws.Range("I122").Value = 2912.0625
$ws.Range("J122").Value = 2179
$ws.Range("K122").Value = 8736.1875
$ws.Range("L122").Value = 6537
$ws.Range("M122").Value = -6286.1875
$ws.Range("N122").Value = -11437

# Row 132: Hull Lotta Damage / Ginseng Lumber
$ws.Range("H132").Value = 32262306
$ws.Range("I132").Value = 40004210
$ws.Range("K132").Value = 120012630
$ws.Range("M132").Value = -120010100

$ws = $wb.Worksheets.Item("CUL")
# Row 7: It's Always Sunny in Vylbrand / Raisins
$ws.Range("H7").Value = 1669405
$ws.Range("I7").Value = 3334329.8
$ws.Range("K7").Value = 10002989.4
$ws.Range("M7").Value = -10002877.4

# Row 23: Sweet Smell of Success / Lavender Oil
$ws.Range("H23").Value = 965.6667
$ws.Range("I23").Value = 900
$ws.Range("K23").Value = 2700
$ws.Range("M23").Value = -2465

# Row 68: Such a Butter Face / Fermented Butter
$ws.Range("H68").Value = 3491.5903
$ws.Range("J68").Value = 3657.234
$ws.Range("L68").Value = 10971.702
$ws.Range("N68").Value = -12593.702

# Row 71: No Margarine of Error (L) / Fermented Butter
$ws.Range("H71").Value = 3491.5903
$ws.Range("J71").Value = 3657.234
$ws.Range("L71").Value = 32915.106
$ws.Range("N71").Value = -41027.106

# Row 122: Salt of the North / Northern Sea Salt
$ws.Range("H122").Value = 893
$ws.Range("I122").Value = 812.1111
$ws.Range("J122").Value = 1014.3333
$ws.Range("K122").Value = 7308.9999
$ws.Range("L122").Value = 9128.9997
$ws.Range("M122").Value = -4858.9999
$ws.Range("N122").Value = -14028.9997

# Row 126: Imperial Palate / Glory Be Soup
$ws.Range("H126").Value = 111111110
$ws.Range("I126").Value = 111111110
$ws.Range("K126").Value = 333333330
$ws.Range("M126").Value = -333328390

# Row 132: More Mezcal / Cooking Mezcal
$ws.Range("H132").Value = 11440.333
$ws.Range("J132").Value = 1000
$ws.Range("L132").Value = 9000
$ws.Range("N132").Value = -14060

$ws = $wb.Worksheets.Item("GSM")
# Row 49: Faith and Fashion / Mythril Earrings
$ws.Range("H49").Value = 30000
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 4: Sole Traders / Leather Duckbills
$ws.Range("H4").Value = 33333
$ws.Range("I4").Value = 33333
$ws.Range("K4").Value = 33333
$ws.Range("M4").Value = -33220

# Row 5: These Boots Are Made for Wailing / Leather Duckbills of Gathering
$ws.Range("H5").Value = 33333
$ws.Range("I5").Value = 33333
$ws.Range("K5").Value = 33333
$ws.Range("M5").Value = -33220

# Row 28: My Sole to Take / Padded Leather Duckbills
$ws.Range("H28").Value = 33333
$ws.Range("I28").Value = 33333
$ws.Range("K28").Value = 33333
$ws.Range("M28").Value = -33101

# Row 37: Quicker than Sand / Padded Leather Duckbills
$ws.Range("H37").Value = 33333
$ws.Range("I37").Value = 33333
$ws.Range("K37").Value = 33333
$ws.Range("M37").Value = -33226

# Row 48: Through a Glass Brightly / Fingerless Peisteskin Gloves
$ws.Range("H48").Value = 33333
$ws.Range("I48").Value = 33333
$ws.Range("K48").Value = 33333
$ws.Range("M48").Value = -32672

# Row 100: Tiger in the Sack / Tiger Leather
$ws.Range("H100").Value = 10938936
$ws.Range("I100").Value = 17501704
$ws.Range("K100").Value = 17501704
$ws.Range("M100").Value = -17501163

# Row 105: Thick and Thin / Gazelleskin Corselet of Scouting
$ws.Range("H105").Value = 40000
$ws.Range("J105").Value = 40000
$ws.Range("L105").Value = 40000
$ws.Range("N105").Value = -46988

# Row 132: Tenets of Tanning / Silver Lobo Leather
$ws.Range("H132").Value = 6107903.5
$ws.Range("I132").Value = 10428540
$ws.Range("J132").Value = 8180.7646
$ws.Range("K132").Value = 31285620
$ws.Range("L132").Value = 24542.2938
$ws.Range("M132").Value = -31283090
$ws.Range("N132").Value = -29602.2938

# Row 136: Respect for Br'aax / Br'aax Leather
$ws.Range("H136").Value = 3625.6843
$ws.Range("I136").Value = 3840.4707
$ws.Range("J136").Value = 1800
$ws.Range("K136").Value = 11521.4121
$ws.Range("L136").Value = 5400
$ws.Range("M136").Value = -8971.4121
$ws.Range("N136").Value = -10500

$ws = $wb.Worksheets.Item("WVR")
# Row 24: Touch Me If You Can / Cotton Work Gloves
$ws.Range("H24").Value = 33333
$ws.Range("I24").Value = 33333
$ws.Range("K24").Value = 33333
$ws.Range("M24").Value = -33103

# Row 26: New Shoes, New Me / Cotton Dress Shoes
$ws.Range("H26").Value = 31110.666
$ws.Range("J26").Value = 29999
$ws.Range("L26").Value = 29999
$ws.Range("N26").Value = -30585

# Row 107: Flax Wax / Bright Linen Yarn
$ws.Range("H107").Value = 420.5
$ws.Range("I107").Value = 391.30768
$ws.Range("K107").Value = 1173.92304
$ws.Range("M107").Value = 746.0769599999999

# Row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 62512200
$ws.Range("I132").Value = 83336510
$ws.Range("J132").Value = 39250
$ws.Range("K132").Value = 250009530
$ws.Range("L132").Value = 117750
$ws.Range("M132").Value = -250007000
$ws.Range("N132").Value = -122810

# Row 136: Weaving the Envelope / Sarcenet Cloth
$ws.Range("H136").Value = 29413608
$ws.Range("J136").Value = 3000
$ws.Range("L136").Value = 9000
$ws.Range("N136").Value = -14100
